$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.6
$ws.Range("G2").Value = 3.9
$ws.Range("H2").Value = 2.34
$ws.Range("I2").Value = 2.44
$ws.Range("J2").Value = 3.05
$ws.Range("K2").Value = 3.2
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 5.8
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 2.02
$ws.Range("Q2").Value = 1.94
$ws.Range("S2").Value = 4.1
$ws.Range("T2").Value = 1.32
$ws.Range("U2").Value = 3.55
$ws.Range("V2").Value = 1.7
$ws.Range("W2").Value = 1.33
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 6.2
$ws.Range("Z2").Value = 16.5
$ws.Range("AA2").Value = 75
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 4.5
$ws.Range("AD2").Value = 10.5
$ws.Range("AE2").Value = 55
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 7.2
$ws.Range("AH2").Value = 14
$ws.Range("AI2").Value = 85
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 27
$ws.Range("AL2").Value = 44
$ws.Range("AM2").Value = 210
$ws.Range("AN2").Value = 48
$ws.Range("AO2").Value = 130
# Row 3
$ws.Range("H3").Value = 4.1
$ws.Range("I3").Value = 4.3
$ws.Range("J3").Value = 3.6
$ws.Range("K3").Value = 3.7
$ws.Range("L3").Value = 1.46
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 3.55
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 1.85
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 1.32
$ws.Range("S3").Value = 3.85
$ws.Range("T3").Value = 1.89
$ws.Range("U3").Value = 2.02
$ws.Range("V3").Value = 1.3
$ws.Range("W3").Value = 1.92
$ws.Range("X3").Value = 14
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 30
$ws.Range("AA3").Value = 140
$ws.Range("AB3").Value = 8.8
$ws.Range("AC3").Value = 7.8
$ws.Range("AD3").Value = 17
$ws.Range("AE3").Value = 70
$ws.Range("AF3").Value = 12
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 19.5
$ws.Range("AI3").Value = 95
$ws.Range("AJ3").Value = 25
$ws.Range("AK3").Value = 23
$ws.Range("AM3").Value = 140
$ws.Range("AN3").Value = 18
# Row 4
$ws.Range("F4").Value = 2.52
$ws.Range("G4").Value = 2.74
$ws.Range("H4").Value = 3.05
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 3.4
$ws.Range("L4").Value = 1.49
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 3.2
$ws.Range("P4").Value = 1.74
$ws.Range("Q4").Value = 2.18
$ws.Range("R4").Value = 1.27
$ws.Range("S4").Value = 4.1
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.42
$ws.Range("W4").Value = 1.58
$ws.Range("X4").Value = 12.5
$ws.Range("Z4").Value = 26
$ws.Range("AA4").Value = 90
$ws.Range("AB4").Value = 9.6
$ws.Range("AC4").Value = 7.2
$ws.Range("AD4").Value = 14.5
$ws.Range("AF4").Value = 17
$ws.Range("AH4").Value = 19.5
$ws.Range("AI4").Value = 65
$ws.Range("AJ4").Value = 50
$ws.Range("AK4").Value = 40
$ws.Range("AN4").Value = 36
$ws.Range("AO4").Value = 48
# Row 5
$ws.Range("F5").Value = 1.34
$ws.Range("G5").Value = 1.35
$ws.Range("H5").Value = 10.5
$ws.Range("I5").Value = 13
$ws.Range("J5").Value = 5.8
$ws.Range("K5").Value = 6.4
$ws.Range("L5").Value = 1.36
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 4.7
$ws.Range("P5").Value = 2.3
$ws.Range("Q5").Value = 1.69
$ws.Range("R5").Value = 1.5
$ws.Range("S5").Value = 2.78
$ws.Range("T5").Value = 2.06
$ws.Range("U5").Value = 1.81
$ws.Range("V5").Value = 1.09
$ws.Range("W5").Value = 3.85
$ws.Range("X5").Value = 22
$ws.Range("Y5").Value = 40
$ws.Range("Z5").Value = 110
$ws.Range("AB5").Value = 9
$ws.Range("AC5").Value = 13
$ws.Range("AD5").Value = 42
$ws.Range("AE5").Value = 200
$ws.Range("AF5").Value = 8.8
$ws.Range("AH5").Value = 32
$ws.Range("AJ5").Value = 11
$ws.Range("AK5").Value = 14.5
$ws.Range("AL5").Value = 38
$ws.Range("AN5").Value = 5.7
# Row 6
$ws.Range("F6").Value = 3.4
$ws.Range("G6").Value = 3.6
$ws.Range("H6").Value = 2.18
$ws.Range("I6").Value = 2.28
$ws.Range("J6").Value = 3.65
$ws.Range("K6").Value = 3.75
$ws.Range("N6").Value = 3.9
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 1.98
$ws.Range("Q6").Value = 1.96
$ws.Range("T6").Value = 1.76
$ws.Range("U6").Value = 2.14
$ws.Range("V6").Value = 1.78
$ws.Range("W6").Value = 1.38
$ws.Range("X6").Value = 15.5
$ws.Range("Y6").Value = 10.5
$ws.Range("Z6").Value = 14
$ws.Range("AA6").Value = 29
$ws.Range("AB6").Value = 14
$ws.Range("AD6").Value = 11
$ws.Range("AE6").Value = 24
$ws.Range("AF6").Value = 26
$ws.Range("AG6").Value = 15
$ws.Range("AJ6").Value = 70
$ws.Range("AK6").Value = 44
$ws.Range("AL6").Value = 55
$ws.Range("AN6").Value = 44
$ws.Range("AO6").Value = 18.5
# Row 7
$ws.Range("F7").Value = 1.28
$ws.Range("G7").Value = 1.29
$ws.Range("H7").Value = 13.5
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 6.2
$ws.Range("K7").Value = 7
$ws.Range("L7").Value = 1.3
$ws.Range("N7").Value = 5.7
$ws.Range("P7").Value = 2.56
$ws.Range("Q7").Value = 1.61
$ws.Range("R7").Value = 1.61
$ws.Range("S7").Value = 2.54
$ws.Range("T7").Value = 2.12
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 1.07
$ws.Range("W7").Value = 4.4
$ws.Range("X7").Value = 26
$ws.Range("Y7").Value = 48
$ws.Range("Z7").Value = 140
$ws.Range("AA7").Value = 690
$ws.Range("AC7").Value = 14.5
$ws.Range("AD7").Value = 46
$ws.Range("AE7").Value = 280
$ws.Range("AF7").Value = 8
$ws.Range("AH7").Value = 36
$ws.Range("AI7").Value = 200
$ws.Range("AJ7").Value = 9.6
$ws.Range("AK7").Value = 14
$ws.Range("AL7").Value = 36
$ws.Range("AM7").Value = 180
$ws.Range("AN7").Value = 4.8
$ws.Range("AO7").Value = 1000
# Row 8
$ws.Range("F8").Value = 1.8
$ws.Range("H8").Value = 5.7
$ws.Range("I8").Value = 6.8
$ws.Range("J8").Value = 3.3
$ws.Range("N8").Value = 2.46
$ws.Range("O8").Value = 1.61
$ws.Range("P8").Value = 1.48
$ws.Range("Q8").Value = 2.92
$ws.Range("T8").Value = 2.4
$ws.Range("U8").Value = 1.61
$ws.Range("X8").Value = 9.6
$ws.Range("Y8").Value = 1000
$ws.Range("AB8").Value = 6
$ws.Range("AC8").Value = 1000
$ws.Range("AF8").Value = 1000
$ws.Range("AG8").Value = 1000
# Row 9
$ws.Range("F9").Value = 2.06
$ws.Range("H9").Value = 3.35
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 3.9
$ws.Range("P9").Value = 2.08
$ws.Range("Q9").Value = 1.78
$ws.Range("S9").Value = 2.96
$ws.Range("T9").Value = 1.65
$ws.Range("U9").Value = 2.12
$ws.Range("V9").Value = 1.35
